$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.1
$ws.Range("H2").Value = 3.2
$ws.Range("I2").Value = 2.05
$ws.Range("J2").Value = 4.75
$ws.Range("AA2").Value = 2.1
$ws.Range("AB2").Value = 1.67
$ws.Range("AC2").Value = 9
$ws.Range("AD2").Value = 19
$ws.Range("AE2").Value = 15
$ws.Range("AG2").Value = 41
$ws.Range("AO2").Value = 8.5
$ws.Range("AQ2").Value = 17
$ws.Range("AR2").Value = 19
$ws.Range("G3").Value = 1.8
$ws.Range("H3").Value = 3.25
$ws.Range("I3").Value = 5.25
$ws.Range("L3").Value = 5.5
$ws.Range("Q3").Value = 1.88
$ws.Range("R3").Value = 1.98
$ws.Range("S3").Value = 2.5
$ws.Range("T3").Value = 1.5
$ws.Range("U3").Value = 3.9
$ws.Range("V3").Value = 1.26
$ws.Range("W3").Value = 5
$ws.Range("X3").Value = 1.17
$ws.Range("AA3").Value = 2.2
$ws.Range("AB3").Value = 1.62
$ws.Range("AC3").Value = 5.5
$ws.Range("AN3").Value = 10
$ws.Range("AO3").Value = 23
$ws.Range("AP3").Value = 17
$ws.Range("AQ3").Value = 51
$ws.Range("AR3").Value = 41
$ws.Range("G5").Value = 2.6
$ws.Range("H5").Value = 3.2
$ws.Range("J5").Value = 3.25
$ws.Range("L5").Value = 3.5
$ws.Range("S5").Value = 2.2
$ws.Range("T5").Value = 1.67
$ws.Range("AA5").Value = 1.83
$ws.Range("AB5").Value = 1.83
$ws.Range("AI5").Value = 8
$ws.Range("AK5").Value = 15
$ws.Range("AM5").Value = 351
$ws.Range("AN5").Value = 8
$ws.Range("AR5").Value = 23
$ws.Range("AS5").Value = 34
$ws.Range("H6").Value = 3.1
$ws.Range("I6").Value = 2.3
$ws.Range("J6").Value = 4
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 3.1
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 8
$ws.Range("O6").Value = 1.4
$ws.Range("P6").Value = 2.75
$ws.Range("S6").Value = 2.25
$ws.Range("T6").Value = 1.62
$ws.Range("W6").Value = 4
$ws.Range("X6").Value = 1.22
$ws.Range("Y6").Value = 1.5
$ws.Range("Z6").Value = 2.5
$ws.Range("AA6").Value = 1.91
$ws.Range("AB6").Value = 1.8
$ws.Range("AC6").Value = 8.5
$ws.Range("AD6").Value = 15
$ws.Range("AG6").Value = 29
$ws.Range("AH6").Value = 41
$ws.Range("AI6").Value = 7.5
$ws.Range("AJ6").Value = 6
$ws.Range("AK6").Value = 15
$ws.Range("AL6").Value = 51
$ws.Range("AM6").Value = 351
$ws.Range("AN6").Value = 7
$ws.Range("AO6").Value = 10
$ws.Range("AP6").Value = 10
$ws.Range("AR6").Value = 21
$ws.Range("AS6").Value = 34
$ws.Range("G7").Value = 3.3
$ws.Range("I7").Value = 2.15
$ws.Range("AA7").Value = 1.73
$ws.Range("AB7").Value = 2
$ws.Range("AM7").Value = 201
$ws.Range("AO7").Value = 10
$ws.Range("AQ7").Value = 19
$ws.Range("AR7").Value = 17
$ws.Range("AS7").Value = 26
$ws.Range("G9").Value = 8.5
$ws.Range("H9").Value = 6.1
$ws.Range("I9").Value = 1.26
$ws.Range("J9").Value = 6.6
$ws.Range("K9").Value = 2.95
$ws.Range("L9").Value = 1.62
$ws.Range("Y9").Value = 1.19
$ws.Range("Z9").Value = 4.15
$ws.Range("AC9").Value = 37
$ws.Range("AD9").Value = 75
$ws.Range("AE9").Value = 27
$ws.Range("AF9").Value = 200
$ws.Range("AG9").Value = 80
$ws.Range("AH9").Value = 55
$ws.Range("AN9").Value = 13
$ws.Range("AO9").Value = 9.25
